{"js": "// Apply the Canada LegendPlay review-brief control-sheet fixes:\n//  1. Replace the ConnexOntario mention in the \"What to include\" bullet list\n//     with pan-Canadian resources wording.\n//  2. Replace the compliance-checklist \"Helpline\" bullet with a\n//     \"Helplines: Provincial helplines (...)\" bullet and add a brand-new\n//     \"National Resource: Responsible Gambling Council\" bullet right after it.\n//  3. Fix the \"parlay calculator\" internal-link target.\n//  4. Fix the \"odds calculator\" internal-link target.\n\nconst body = context.document.body;\n\n// --- helper: replace the first search hit's text (exact, case-sensitive) ---\nasync function replaceFirst(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Canadian responsible gambling resources bullet (bulleted \"include\" list).\nawait replaceFirst(\n  \"Canadian responsible gambling resources (ConnexOntario: 1-800-463-1554)\",\n  \"Canadian responsible gambling resources (Provincial helplines and Responsible Gambling Council)\"\n);\n\n// 2) Compliance checklist \"Helpline\" bullet -> \"Helplines\" bullet, plus a new\n//    bullet paragraph for the national resource right after it.\nconst helplineHits = body.search(\"[ ] Helpline: 1-800-463-1554 (ConnexOntario)\", { matchCase: true });\nhelplineHits.load(\"items\");\nawait context.sync();\nif (helplineHits.items.length === 0) {\n  throw new Error(\"Helpline bullet not found\");\n}\nconst helplineRange = helplineHits.items[0];\nconst helplineParagraph = helplineRange.paragraphs.getFirst();\nhelplineRange.insertText(\n  \"[ ] Helplines: Provincial helplines (AB: 1-866-332-2322, BC: 1-888-795-6111, QC: 1-800-461-0140)\",\n  \"Replace\"\n);\nawait context.sync();\n\nhelplineParagraph.insertParagraph(\n  \"[ ] National Resource: Responsible Gambling Council (www.responsiblegambling.org)\",\n  \"After\"\n);\nawait context.sync();\n\n// 3) \"parlay calculator\" internal-link target.\nawait replaceFirst(\n  '\"parlay calculator\" \\u2192 /sport/betting/calculators/parlay.htm',\n  '\"parlay calculator\" \\u2192 /sport/betting-tools/parlay-calculator.htm.htm'\n);\n\n// 4) \"odds calculator\" internal-link target.\nawait replaceFirst(\n  '\"odds calculator\" \\u2192 /sport/betting/calculators/odds.htm',\n  '\"odds calculator\" \\u2192 /sport/betting-tools/odds-calculator.htm.htm'\n);\n", "ps1": "# Apply the Canada LegendPlay review-brief control-sheet fixes:\n#  1. Replace the ConnexOntario mention in the \"what to include\" bullet list\n#     with pan-Canadian resources wording.\n#  2. Replace the compliance-checklist \"Helpline\" bullet with a\n#     \"Helplines: Provincial helplines (...)\" bullet and add a brand-new\n#     \"National Resource: Responsible Gambling Council\" bullet right after it.\n#  3. Fix the \"parlay calculator\" internal-link target.\n#  4. Fix the \"odds calculator\" internal-link target.\n#\n# NOTE: text is written directly onto Paragraph.Range.Text (not via\n# Find.Replacement, which runs the text through AutoCorrect's \"smart quotes\"\n# substitution and would mangle the straight double-quotes in items 3/4), and\n# paragraphs are matched with String.Contains (not the -like operator, whose\n# \"[ ]\" wildcard character class would otherwise swallow the literal\n# \"[ ] \" checkbox prefix used throughout this checklist).\n\n$d = $word.ActiveDocument\n$arrow = [char]0x2192\n\n# 1) Canadian responsible gambling resources bullet (bulleted \"include\" list).\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains(\"Canadian responsible gambling resources (ConnexOntario: 1-800-463-1554)\")) {\n        $p.Range.Text = \"Canadian responsible gambling resources (Provincial helplines and Responsible Gambling Council)\"\n        break\n    }\n}\n\n# 2) Compliance checklist \"Helpline\" bullet -> \"Helplines\" bullet, plus a new\n#    bullet paragraph for the national resource right after it.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains(\"[ ] Helpline: 1-800-463-1554 (ConnexOntario)\")) {\n        $p.Range.Text = \"[ ] Helplines: Provincial helplines (AB: 1-866-332-2322, BC: 1-888-795-6111, QC: 1-800-461-0140)\"\n        $p.Range.InsertParagraphAfter() | Out-Null\n        $p.Next().Range.Text = \"[ ] National Resource: Responsible Gambling Council (www.responsiblegambling.org)\"\n        break\n    }\n}\n\n# 3) \"parlay calculator\" internal-link target.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains('\"parlay calculator\" ' + $arrow + ' /sport/betting/calculators/parlay.htm')) {\n        $p.Range.Text = '\"parlay calculator\" ' + $arrow + ' /sport/betting-tools/parlay-calculator.htm.htm'\n        break\n    }\n}\n\n# 4) \"odds calculator\" internal-link target.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains('\"odds calculator\" ' + $arrow + ' /sport/betting/calculators/odds.htm')) {\n        $p.Range.Text = '\"odds calculator\" ' + $arrow + ' /sport/betting-tools/odds-calculator.htm.htm'\n        break\n    }\n}\n"}
